# week12 notes on temporary diff using q
# Applies the edits made on the "Dierencia Temporal sin fπ(s)" sheet:
#  - drops the scratch "s→sf" label/highlight in Q6
#  - updates the Q8 comparison row (and N7) with new numbers
#  - removes the leftover scratch formula in N11
#  - leaves the selection where the author ended up (K13, scrolled to show column E onward)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dierencia Temporal sin fπ(s)")
$ws.Activate()

# Remove the "s→sf" label together with its highlight fill in Q6
$ws.Range("Q6").Clear()

# New SARSA / Q-learning comparison values
$ws.Range("N7").Value = 8.9969999999999999

$ws.Range("M8").Value = 0.67
$ws.Range("N8").Value = 0.15
$ws.Range("O8").Value = 0.8
$ws.Range("P8").Value = 0.53500000000000003
$ws.Range("Q8").Value = 0.89

# Drop the scratch formula that had been left in N11
$ws.Range("N11").ClearContents()

# Restore the view/selection to where the author left it
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("K13").Select() | Out-Null
